# Updated Extended Attributes list
# Added "Upload file name" because it exists in the Test environment

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtendedAttributes")

# Add the new row describing the "Upload File Name" extended attribute
$ws.Range("A17").Value = "NA"
$ws.Range("B17").Value = "Upload File Name"
$ws.Range("C17").Value = "TEXT"
$ws.Range("D17").Value = "OBSERVATION"
$ws.Range("E17").Value = "NA"
$ws.Range("F17").Value = "NA"
$ws.Range("G17").Value = "NA"
$ws.Range("H17").Value = "NA"
$ws.Range("I17").Value = $false
$ws.Range("J17").Value = "Used by EDT to tag observations by what upload file they came from."

# Update the (previously active) DropdownLists sheet's saved selection
# first, then finish on ExtendedAttributes so it ends up as the active tab
# (mirrors the author's saved view state).
$dropdownWs = $wb.Worksheets.Item("DropdownLists")
$dropdownWs.Range("A18").Select()

$ws.Activate()
$ws.Range("A18:XFD1048576").Select()
